{"js": "// Office.js (Word JavaScript API) script.\n//\n// The underlying author edit (per the commit's XML diff) boils down to two\n// real textual changes inside the report body; the remaining hunks in the\n// diff (title run, \"buffer\" run, header run) are just Word's own artifact\n// of merging adjacent same-formatted runs / dropping stale <w:proofErr>\n// spell-check bookmarks around the edited text \u2014 the rendered text in\n// those spots is byte-for-byte identical before and after, so nothing\n// needs to change there.\n//\n//  1) \"...je ne suis pas parvenu \u00e0 r\u00e9 utiliser fetch.\" -> \"...parvenue...\"\n//     (grammar/agreement fix: insert \"e\" after \"parvenu\").\n//  2) New sentence inserted into the \"l\u00e9gende\" paragraph, between the\n//     \"...funiculaires.\" sentence and the \"J'ai pens\u00e9...\" sentence:\n//     \" Toujours quand on clique sur les lignes transports, un encadr\u00e9\n//     bleu appara\u00eet pour mettre en avant la ligne de transport choisi.\"\n\nconst body = context.document.body;\n\n// --- 1) \"parvenu\" -> \"parvenue\" (gender agreement) ------------------------\nconst parvenuHits = body.search(\"parvenu \u00e0 r\u00e9 utiliser fetch.\", { matchCase: true });\nparvenuHits.load(\"text\");\nawait context.sync();\n\nif (parvenuHits.items.length > 0) {\n  parvenuHits.items[0].insertText(\"parvenue \u00e0 r\u00e9 utiliser fetch.\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2) Insert new sentence about the blue highlight box ------------------\nconst legendeHits = body.search(\n  \"funiculaires. J\\u2019ai pens\\u00e9 \\u00e0 en faire une pour les \\u00e9quipements sportifs\",\n  { matchCase: true }\n);\nlegendeHits.load(\"text\");\nawait context.sync();\n\nif (legendeHits.items.length > 0) {\n  const replacement =\n    \"funiculaires.\" +\n    \" Toujours quand on clique sur les lignes transports, un encadr\u00e9 bleu appara\u00eet pour mettre en avant la ligne de transport choisi.\" +\n    \" J\\u2019ai pens\u00e9 \u00e0 en faire une pour les \u00e9quipements sportifs\";\n  legendeHits.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# The underlying author edit (per the commit's XML diff) boils down to two\n# real textual changes inside the report body; the remaining hunks in the\n# diff (title run, \"buffer\" run, header run) are just Word's own artifact\n# of merging adjacent same-formatted runs / dropping stale proofing-error\n# bookmarks around the edited text -- the rendered text in those spots is\n# identical before and after, so nothing needs to change there.\n#\n#  1) \"...je ne suis pas parvenu a re utiliser fetch.\" -> \"...parvenue...\"\n#     (grammar/agreement fix: insert \"e\" after \"parvenu\").\n#  2) New sentence inserted into the \"legende\" paragraph, between the\n#     \"...funiculaires.\" sentence and the \"J'ai pense...\" sentence.\n\n$d = $word.ActiveDocument\n\n# --- 1) \"parvenu\" -> \"parvenue\" (gender agreement) -------------------------\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.MatchCase = $true\n$find1.Find.MatchWholeWord = $false\n$find1.Find.MatchWildcards = $false\n$find1.Find.Execute(\"parvenu \u00e0 r\u00e9 utiliser fetch.\", $false, $true, $false, $false, $false, $true, 1, $false, \"parvenue \u00e0 r\u00e9 utiliser fetch.\", 2)\n\n# --- 2) Insert new sentence about the blue highlight box --------------------\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.MatchCase = $true\n$find2.Find.MatchWholeWord = $false\n$find2.Find.MatchWildcards = $false\n$find2.Find.Execute(\"et funiculaires. J\u2019ai pens\u00e9\", $false, $true, $false, $false, $false, $true, 1, $false, \"et funiculaires. Toujours quand on clique sur les lignes transports, un encadr\u00e9 bleu appara\u00eet pour mettre en avant la ligne de transport choisi. J\u2019ai pens\u00e9\", 2)\n\n$d.Save()\n"}
